# Applies the cryptos list refresh (updated prices/volumes and a few
# re-ranked rows) as captured in the authoritative XML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @('D2', '66.492.45'),
  @('E2', '  +0.49%  '),
  @('D3', '3.608.86'),
  @('E3', '  +1.19%  '),
  @('E4', '  -0.02%  '),
  @('D5', '608.14'),
  @('E5', '  +0.38%  '),
  @('D6', '148.99'),
  @('E6', '  +2.96%  '),
  @('D7', '3.605.72'),
  @('E7', '  +1.17%  '),
  @('E8', '  +0.02%  '),
  @('D9', '0.487'),
  @('E9', '  -0.46%  '),
  @('D10', '0.136'),
  @('E10', '  +0.32%  '),
  @('D11', '7.99'),
  @('E11', '  +1.59%  '),
  @('D12', '0.415'),
  @('E12', '  +0.51%  '),
  @('D13', '4.219.61'),
  @('E13', '  +1.27%  '),
  @('D14', '0.0000207'),
  @('E14', '  +0.15%  '),
  @('D15', '29.83'),
  @('E15', '  -0.74%  '),
  @('D16', '3.616.94'),
  @('E16', '  +1.46%  '),
  @('E17', '  +2.16%  '),
  @('D18', '66.536.79'),
  @('E18', '  +0.45%  '),
  @('D19', '11.29'),
  @('E19', '  -2.14%  '),
  @('D20', '6.35'),
  @('E20', '  +2.07%  '),
  @('D21', '14.98'),
  @('E21', '  +1.70%  '),
  @('D22', '425.00'),
  @('E22', '  -1.18%  '),
  @('D23', '0.614'),
  @('E23', '  +0.47%  '),
  @('D24', '78.65'),
  @('E24', '  -0.76%  '),
  @('E25', '  -0.05%  '),
  @('D26', '0.0000122'),
  @('E26', '  +3.84%  '),
  @('D27', '8.31'),
  @('E27', '  +5.56%  '),
  @('D28', '9.43'),
  @('E28', '  +3.62%  '),
  @('D29', '2.50'),
  @('E29', '  +0.11%  '),
  @('E30', '  -0.04%  '),
  @('B31', 'RenzoRestakedETH'),
  @('C31', 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'),
  @('D31', '3.604.06'),
  @('E31', '  +1.20%  '),
  @('B32', 'Kaspa'),
  @('C32', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'),
  @('D32', '0.159'),
  @('E32', '  +4.14%  '),
  @('B33', 'Fetch.AI'),
  @('C33', 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'),
  @('D33', '1.46'),
  @('E33', '  +1.14%  '),
  @('D34', '25.27'),
  @('E34', '  -0.86%  '),
  @('B35', 'USDe'),
  @('C35', 'https://coinranking.com/coin/exbfr2U-0+usde-usde'),
  @('D35', '1.00'),
  @('E35', '  +0.00%  '),
  @('B36', 'Aptos'),
  @('C36', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'),
  @('D36', '7.82'),
  @('E36', '  -0.06%  '),
  @('D37', '5.64'),
  @('E37', '  +1.48%  '),
  @('D38', '1.69'),
  @('E38', '  -2.20%  '),
  @('D39', '176.58'),
  @('E39', '  +1.60%  '),
  @('D40', '0.0857'),
  @('E40', '  +1.52%  '),
  @('D41', '5.21'),
  @('E41', '  +0.47%  '),
  @('D42', '0.893'),
  @('E42', '  -0.33%  '),
  @('B43', 'OKB'),
  @('C43', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'),
  @('D43', '45.94'),
  @('E43', '  -0.02%  '),
  @('B44', 'Stacks'),
  @('C44', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'),
  @('D44', '1.87'),
  @('E44', '  -4.40%  '),
  @('D45', '2.57'),
  @('E45', '  +6.72%  '),
  @('D46', '0.999'),
  @('E46', '  +0.00%  '),
  @('D47', '24.58'),
  @('E47', '  -1.71%  '),
  @('D48', '23.73'),
  @('E48', '  +3.54%  '),
  @('B49', 'Cosmos'),
  @('C49', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'),
  @('D49', '7.18'),
  @('E49', '  +0.51%  '),
  @('B50', 'ONDO'),
  @('C50', 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'),
  @('D50', '1.14'),
  @('E50', '  -4.80%  '),
  @('B51', 'SuiNetwork'),
  @('C51', 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'),
  @('D51', '0.969'),
  @('E51', '  +4.04%  ')
)

foreach ($item in $updates) {
  $cellRef = $item[0]
  $newVal = $item[1]
  $rng = $ws.Range($cellRef)
  if ($cellRef.StartsWith("D")) {
    # Price column: values like "425.00" or "66.492.45" must stay text,
    # otherwise Excel auto-coerces them to numbers and drops formatting.
    $rng.NumberFormat = "@"
    $rng.Value = $newVal
    $rng.Style = "Normal"
  } else {
    $rng.Value = $newVal
  }
}
